$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("reviews_count"), shifting columns F:K left to E:J
$ws.Range("E:E").Delete()
